$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two records were dropped from the dataset entirely: "RM 232" (row 26) and
# "SC 92" (row 28). Delete the lower-numbered row second so its index is
# still valid when we get to it.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# The remaining rows shift up to fill rows 26-33. Refresh the "missing
# value" mask for column D (sheet column E) and column A (sheet column B):
# some previously-missing values now have numbers, and some previously
# present values are now missing.
$ws.Range("E2").Value = -7.2
$ws.Range("E6").Value = ""
$ws.Range("E12").Value = -5.3
$ws.Range("E14").Value = ""
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.7
$ws.Range("E23").Value = ""
$ws.Range("E24").Value = ""

$ws.Range("B26").Value = -20.2
$ws.Range("B27").Value = ""
$ws.Range("B28").Value = ""
$ws.Range("B29").Value = -19.5
$ws.Range("B30").Value = -19.7
$ws.Range("B31").Value = ""
$ws.Range("E31").Value = -8.1
$ws.Range("B32").Value = ""
$ws.Range("E33").Value = -10.7
